$d = $word.ActiveDocument

$d.Content.Find.Execute("72÷3=24, 0", $true, $false, $false, $false, $false, $true, 1, $false, "29÷3=9, 2", 2) | Out-Null
$d.Content.Find.Execute("20÷6=3, 2", $true, $false, $false, $false, $false, $true, 1, $false, "81÷6=13, 3", 2) | Out-Null
$d.Content.Find.Execute("42÷5=8, 2", $true, $false, $false, $false, $false, $true, 1, $false, "15÷5=3, 0", 2) | Out-Null
$d.Content.Find.Execute("55÷9=6, 1", $true, $false, $false, $false, $false, $true, 1, $false, "80÷9=8, 8", 2) | Out-Null
$d.Content.Find.Execute("78÷2=39, 0", $true, $false, $false, $false, $false, $true, 1, $false, "21÷9=2, 3", 2) | Out-Null
$d.Content.Find.Execute("96÷4=24, 0", $true, $false, $false, $false, $false, $true, 1, $false, "46÷2=23, 0", 2) | Out-Null
$d.Content.Find.Execute("60÷7=8, 4", $true, $false, $false, $false, $false, $true, 1, $false, "27÷4=6, 3", 2) | Out-Null
$d.Content.Find.Execute("57÷5=11, 2", $true, $false, $false, $false, $false, $true, 1, $false, "18÷7=2, 4", 2) | Out-Null
$d.Content.Find.Execute("26÷6=4, 2", $true, $false, $false, $false, $false, $true, 1, $false, "82÷9=9, 1", 2) | Out-Null
$d.Content.Find.Execute("78÷8=9, 6", $true, $false, $false, $false, $false, $true, 1, $false, "69÷5=13, 4", 2) | Out-Null
$d.Content.Find.Execute("87÷2=43, 1", $true, $false, $false, $false, $false, $true, 1, $false, "40÷9=4, 4", 2) | Out-Null
$d.Content.Find.Execute("84÷9=9, 3", $true, $false, $false, $false, $false, $true, 1, $false, "51÷6=8, 3", 2) | Out-Null
$d.Content.Find.Execute("92÷4=23, 0", $true, $false, $false, $false, $false, $true, 1, $false, "63÷4=15, 3", 2) | Out-Null
$d.Content.Find.Execute("13÷5=2, 3", $true, $false, $false, $false, $false, $true, 1, $false, "21÷9=2, 3", 2) | Out-Null
$d.Content.Find.Execute("79÷4=19, 3", $true, $false, $false, $false, $false, $true, 1, $false, "61÷2=30, 1", 2) | Out-Null
$d.Content.Find.Execute("44÷9=4, 8", $true, $false, $false, $false, $false, $true, 1, $false, "57÷5=11, 2", 2) | Out-Null
$d.Content.Find.Execute("57÷9=6, 3", $true, $false, $false, $false, $false, $true, 1, $false, "88÷4=22, 0", 2) | Out-Null
$d.Content.Find.Execute("52÷7=7, 3", $true, $false, $false, $false, $false, $true, 1, $false, "91÷7=13, 0", 2) | Out-Null
$d.Content.Find.Execute("41÷7=5, 6", $true, $false, $false, $false, $false, $true, 1, $false, "65÷6=10, 5", 2) | Out-Null
$d.Content.Find.Execute("93÷5=18, 3", $true, $false, $false, $false, $false, $true, 1, $false, "48÷7=6, 6", 2) | Out-Null
$d.Content.Find.Execute("42÷8=5, 2", $true, $false, $false, $false, $false, $true, 1, $false, "16÷9=1, 7", 2) | Out-Null
$d.Content.Find.Execute("21÷7=3, 0", $true, $false, $false, $false, $false, $true, 1, $false, "42÷2=21, 0", 2) | Out-Null
$d.Content.Find.Execute("63÷9=7, 0", $true, $false, $false, $false, $false, $true, 1, $false, "27÷8=3, 3", 2) | Out-Null
$d.Content.Find.Execute("44÷5=8, 4", $true, $false, $false, $false, $false, $true, 1, $false, "50÷3=16, 2", 2) | Out-Null
$d.Content.Find.Execute("95÷7=13, 4", $true, $false, $false, $false, $false, $true, 1, $false, "75÷6=12, 3", 2) | Out-Null
